# Updated cryptos list on Sat Sep  2 13:57:45 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct value updates (text/link/percentage cells) ---
# Row 2
$ws.Cells.Item(2, 5).Value = '  -0.97%  '
# Row 3
$ws.Cells.Item(3, 5).Value = '  -0.67%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.81%  '
# Row 5
$ws.Cells.Item(5, 5).Value = '  -0.03%  '
# Row 6
$ws.Cells.Item(6, 5).Value = '  -1.64%  '
# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.50%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.52%  '
# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.61%  '
# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.93%  '
# Row 11
$ws.Cells.Item(11, 5).Value = '  -1.08%  '
# Row 12
$ws.Cells.Item(12, 2).Value = 'Polkadot'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(12, 5).Value = '  -0.47%  '
# Row 13
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 5).Value = '  -0.28%  '
# Row 14
$ws.Cells.Item(14, 5).Value = '  -0.81%  '
# Row 15
$ws.Cells.Item(15, 5).Value = '  -0.57%  '
# Row 16
$ws.Cells.Item(16, 4).Value = '0.0₅7906'
$ws.Cells.Item(16, 5).Value = '  -1.22%  '
# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.42%  '
# Row 18
$ws.Cells.Item(18, 5).Value = '  -0.91%  '
# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.46%  '
# Row 20
$ws.Cells.Item(20, 5).Value = '  -2.99%  '
# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.17%  '
# Row 22
$ws.Cells.Item(22, 5).Value = '  -1.69%  '
# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.48%  '
# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.40%  '
# Row 25
$ws.Cells.Item(25, 5).Value = '  +5.56%  '
# Row 26
$ws.Cells.Item(26, 5).Value = '  -2.52%  '
# Row 27
$ws.Cells.Item(27, 5).Value = '  -3.03%  '
# Row 28
$ws.Cells.Item(28, 2).Value = 'Cosmos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(28, 5).Value = '  -2.60%  '
# Row 29
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(29, 5).Value = '  -0.87%  '
# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.30%  '
# Row 31
$ws.Cells.Item(31, 5).Value = '  -2.06%  '
# Row 32
$ws.Cells.Item(32, 5).Value = '  -2.12%  '
# Row 33
$ws.Cells.Item(33, 5).Value = '  -1.00%  '
# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.40%  '
# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.64%  '
# Row 36
$ws.Cells.Item(36, 5).Value = '  -3.87%  '
# Row 37
$ws.Cells.Item(37, 5).Value = '  -2.72%  '
# Row 38
$ws.Cells.Item(38, 2).Value = 'Maker'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(38, 5).Value = '  -1.20%  '
# Row 39
$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 5).Value = '  -1.22%  '
# Row 40
$ws.Cells.Item(40, 5).Value = '  -1.09%  '
# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.31%  '
# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.62%  '
# Row 43
$ws.Cells.Item(43, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(43, 5).Value = '  -2.14%  '
# Row 44
$ws.Cells.Item(44, 2).Value = 'Quant'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(44, 5).Value = '  -0.04%  '
# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.02%  '
# Row 46
$ws.Cells.Item(46, 5).Value = '  +4.13%  '
# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.40%  '
# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.11%  '
# Row 49
$ws.Cells.Item(49, 5).Value = '  -0.68%  '
# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.51%  '
# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.43%  '

# --- Price cells that look numeric: force literal text storage like the source ---
# (write as a text formula, then paste-special as values so no formula or number
#  formatting style is left behind on the cell)
$ws.Cells.Item(2, 4).Formula = '="25.879.50"'
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$ws.Cells.Item(3, 4).Formula = '="1.640.96"'
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4163)
$ws.Cells.Item(4, 4).Formula = '="1.000"'
$ws.Cells.Item(4, 4).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4163)
$ws.Cells.Item(5, 4).Formula = '="216.05"'
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(6, 4).Formula = '="0.5040"'
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$ws.Cells.Item(7, 4).Formula = '="1.003"'
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4163)
$ws.Cells.Item(8, 4).Formula = '="0.2578"'
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$ws.Cells.Item(9, 4).Formula = '="0.06390"'
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$ws.Cells.Item(10, 4).Formula = '="19.58"'
$ws.Cells.Item(10, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$ws.Cells.Item(11, 4).Formula = '="0.07752"'
$ws.Cells.Item(11, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 4).Formula = '="4.268"'
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(13, 4).Formula = '="1.647.03"'
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 4).Formula = '="1.864.41"'
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 4).Formula = '="0.5469"'
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 4).Formula = '="64.26"'
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(18, 4).Formula = '="25.911.57"'
$ws.Cells.Item(18, 4).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$ws.Cells.Item(19, 4).Formula = '="1.004"'
$ws.Cells.Item(19, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 4).Formula = '="202.70"'
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 4).Formula = '="4.409"'
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(22, 4).Formula = '="9.898"'
$ws.Cells.Item(22, 4).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$ws.Cells.Item(23, 4).Formula = '="5.980"'
$ws.Cells.Item(23, 4).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4163)
$ws.Cells.Item(24, 4).Formula = '="1.006"'
$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$ws.Cells.Item(25, 4).Formula = '="1.884"'
$ws.Cells.Item(25, 4).Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$ws.Cells.Item(26, 4).Formula = '="141.07"'
$ws.Cells.Item(26, 4).Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$ws.Cells.Item(28, 4).Formula = '="6.798"'
$ws.Cells.Item(28, 4).Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4163)
$ws.Cells.Item(29, 4).Formula = '="15.66"'
$ws.Cells.Item(29, 4).Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4163)
$ws.Cells.Item(30, 4).Formula = '="1.246"'
$ws.Cells.Item(30, 4).Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$ws.Cells.Item(31, 4).Formula = '="0.04977"'
$ws.Cells.Item(31, 4).Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4163)
$ws.Cells.Item(32, 4).Formula = '="3.277"'
$ws.Cells.Item(32, 4).Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4163)
$ws.Cells.Item(33, 4).Formula = '="3.199"'
$ws.Cells.Item(33, 4).Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$ws.Cells.Item(34, 4).Formula = '="1.547"'
$ws.Cells.Item(34, 4).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4163)
$ws.Cells.Item(35, 4).Formula = '="2.373"'
$ws.Cells.Item(35, 4).Copy()
$ws.Cells.Item(35, 4).PasteSpecial(-4163)
$ws.Cells.Item(37, 4).Formula = '="0.8927"'
$ws.Cells.Item(37, 4).Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4163)
$ws.Cells.Item(38, 4).Formula = '="1.152.64"'
$ws.Cells.Item(38, 4).Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4163)
$ws.Cells.Item(39, 4).Formula = '="0.5618"'
$ws.Cells.Item(39, 4).Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4163)
$ws.Cells.Item(40, 4).Formula = '="0.01566"'
$ws.Cells.Item(40, 4).Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$ws.Cells.Item(41, 4).Formula = '="1.006"'
$ws.Cells.Item(41, 4).Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$ws.Cells.Item(42, 4).Formula = '="5.699"'
$ws.Cells.Item(42, 4).Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$ws.Cells.Item(43, 4).Formula = '="0.8091"'
$ws.Cells.Item(43, 4).Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4163)
$ws.Cells.Item(44, 4).Formula = '="100.10"'
$ws.Cells.Item(44, 4).Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$ws.Cells.Item(45, 4).Formula = '="1.775.29"'
$ws.Cells.Item(45, 4).Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4163)
$ws.Cells.Item(47, 4).Formula = '="0.4526"'
$ws.Cells.Item(47, 4).Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$ws.Cells.Item(48, 4).Formula = '="1.007"'
$ws.Cells.Item(48, 4).Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4163)
$ws.Cells.Item(49, 4).Formula = '="54.89"'
$ws.Cells.Item(49, 4).Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4163)
$ws.Cells.Item(50, 4).Formula = '="0.05049"'
$ws.Cells.Item(50, 4).Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$ws.Cells.Item(51, 4).Formula = '="1.003"'
$ws.Cells.Item(51, 4).Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4163)

$excel.CutCopyMode = 0

